$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2892.875
$ws.Range("I113").Value = 2974
$ws.Range("J113").Value = 2865.8333
$ws.Range("K113").Value = 2974
$ws.Range("L113").Value = 2865.8333
$ws.Range("M113").Value = 280
$ws.Range("N113").Value = -9373.8333
$ws.Range("H137").Value = 1701.5
$ws.Range("I137").Value = 1473.3572
$ws.Range("K137").Value = 4420.071599999999
$ws.Range("M137").Value = -1870.071599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3112.9375
$ws.Range("J88").Value = 3152.5
$ws.Range("L88").Value = 3152.5
$ws.Range("N88").Value = -3964.5
$ws.Range("H91").Value = 3112.9375
$ws.Range("J91").Value = 3152.5
$ws.Range("L91").Value = 3152.5
$ws.Range("N91").Value = -5960.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 8750
$ws.Range("I75").Value = 8750
$ws.Range("K75").Value = 8750
$ws.Range("M75").Value = -7814
$ws.Range("H78").Value = 8750
$ws.Range("I78").Value = 8750
$ws.Range("K78").Value = 26250
$ws.Range("M78").Value = -21570
$ws.Range("H80").Value = 410.22223
$ws.Range("I80").Value = 364.2857
$ws.Range("J80").Value = 439.45456
$ws.Range("K80").Value = 364.2857
$ws.Range("L80").Value = 439.45456
$ws.Range("M80").Value = 633.7143
$ws.Range("N80").Value = -2435.45456
$ws.Range("H82").Value = 38000
$ws.Range("I82").Value = 16000
$ws.Range("K82").Value = 16000
$ws.Range("M82").Value = -15617
$ws.Range("H83").Value = 410.22223
$ws.Range("I83").Value = 364.2857
$ws.Range("J83").Value = 439.45456
$ws.Range("K83").Value = 1821.4285
$ws.Range("L83").Value = 2197.2728
$ws.Range("M83").Value = 3170.5715
$ws.Range("N83").Value = -12181.2728
$ws.Range("H85").Value = 38000
$ws.Range("I85").Value = 16000
$ws.Range("K85").Value = 16000
$ws.Range("M85").Value = -14674
$ws.Range("H97").Value = 23137.5
$ws.Range("I97").Value = 22933.334
$ws.Range("K97").Value = 22933.334
$ws.Range("M97").Value = -21942.334
$ws.Range("H141").Value = 98997
$ws.Range("J141").Value = 98997
$ws.Range("L141").Value = 98997
$ws.Range("N141").Value = -109357

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3872.913
$ws.Range("I99").Value = 3258.2
$ws.Range("J99").Value = 5025.5
$ws.Range("K99").Value = 3258.2
$ws.Range("L99").Value = 5025.5
$ws.Range("M99").Value = -1760.2
$ws.Range("N99").Value = -8021.5
$ws.Range("H117").Value = 83000
$ws.Range("J117").Value = 83000
$ws.Range("L117").Value = 83000
$ws.Range("N117").Value = -92178
$ws.Range("H126").Value = 3872.913
$ws.Range("I126").Value = 3258.2
$ws.Range("J126").Value = 5025.5
$ws.Range("K126").Value = 9774.599999999999
$ws.Range("L126").Value = 15076.5
$ws.Range("M126").Value = -7304.599999999999
$ws.Range("N126").Value = -20016.5
$ws.Range("H132").Value = 4004.0435
$ws.Range("I132").Value = 4233.048
$ws.Range("J132").Value = 1599.5
$ws.Range("K132").Value = 12699.144
$ws.Range("L132").Value = 4798.5
$ws.Range("M132").Value = -10169.144
$ws.Range("N132").Value = -9858.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 899.6
$ws.Range("I18").Value = 874.5
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 2623.5
$ws.Range("L18").Value = 3000
$ws.Range("M18").Value = -2454.5
$ws.Range("N18").Value = -3338
$ws.Range("H48").Value = 120
$ws.Range("I48").Value = 120
$ws.Range("K48").Value = 360
$ws.Range("M48").Value = -110
$ws.Range("H113").Value = 209.375
$ws.Range("J113").Value = 219.6
$ws.Range("L113").Value = 658.8
$ws.Range("N113").Value = -4998.8
$ws.Range("H122").Value = 1975.5
$ws.Range("J122").Value = 1975.5
$ws.Range("L122").Value = 17779.5
$ws.Range("N122").Value = -22679.5
$ws.Range("H131").Value = 889.65
$ws.Range("J131").Value = 958.7059
$ws.Range("L131").Value = 2876.1177
$ws.Range("N131").Value = -12956.1177
$ws.Range("H134").Value = 797
$ws.Range("I134").Value = 797
$ws.Range("K134").Value = 2391
$ws.Range("M134").Value = 2679

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H94").Value = 29994
$ws.Range("J94").Value = 29994
$ws.Range("L94").Value = 29994
$ws.Range("N94").Value = -31346
$ws.Range("H107").Value = 784.625
$ws.Range("I107").Value = 579.5
$ws.Range("K107").Value = 579.5
$ws.Range("M107").Value = 1340.5
$ws.Range("H133").Value = 49999
$ws.Range("J133").Value = 49999
$ws.Range("L133").Value = 49999
$ws.Range("N133").Value = -60119
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1416.6666
$ws.Range("I46").Value = 1625
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 1625
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -1437
$ws.Range("N46").Value = -1376
$ws.Range("H61").Value = 2784.7144
$ws.Range("I61").Value = 2784.7144
$ws.Range("K61").Value = 2784.7144
$ws.Range("M61").Value = -2582.7144
$ws.Range("H82").Value = 2932.8667
$ws.Range("I82").Value = 1123.5
$ws.Range("J82").Value = 3590.818
$ws.Range("K82").Value = 1123.5
$ws.Range("L82").Value = 3590.818
$ws.Range("M82").Value = -762.5
$ws.Range("N82").Value = -4312.818
$ws.Range("H85").Value = 2932.8667
$ws.Range("I85").Value = 1123.5
$ws.Range("J85").Value = 3590.818
$ws.Range("K85").Value = 1123.5
$ws.Range("L85").Value = 3590.818
$ws.Range("M85").Value = 124.5
$ws.Range("N85").Value = -6086.818
$ws.Range("H100").Value = 1742.4166
$ws.Range("I100").Value = 1742.4166
$ws.Range("K100").Value = 1742.4166
$ws.Range("M100").Value = -1201.4166
$ws.Range("H113").Value = 2784.7144
$ws.Range("I113").Value = 2784.7144
$ws.Range("K113").Value = 2784.7144
$ws.Range("M113").Value = -614.7143999999998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2339.3
$ws.Range("I96").Value = 1674.375
$ws.Range("J96").Value = 4999
$ws.Range("K96").Value = 1674.375
$ws.Range("L96").Value = 4999
$ws.Range("M96").Value = -301.375
$ws.Range("N96").Value = -7745
$ws.Range("H107").Value = 555.5454999999999
$ws.Range("I107").Value = 442.7143
$ws.Range("K107").Value = 1328.1429
$ws.Range("M107").Value = 591.8571000000002
$ws.Range("H126").Value = 1085.625
$ws.Range("I126").Value = 1085.625
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 3256.875
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -786.875
$ws.Range("N126").ClearContents()
